$wb = $excel.ActiveWorkbook

# --- media sheet: switch the "original" links to the new IIIF Image API ---
$ws = $wb.Worksheets.Item("media")

# Set B3 first, then B2, so the shared-string table picks up the
# "letter2" url at the lower index and "letter1" url at the next index
# (matches the target document's shared-string order).
$ws.Range("B3").Value = "https://05r4t6462c.execute-api.us-east-1.amazonaws.com/latest/iiif/2/tei-eaj%2Fletter2/info.json"
$ws.Range("B2").Value = "https://05r4t6462c.execute-api.us-east-1.amazonaws.com/latest/iiif/2/tei-eaj%2Fletter1/info.json"

# Stash original cell formatting (hyperlink style) so it can be restored
# after Hyperlinks.Add() re-applies its own default hyperlink format.
$ws.Range("B2").Copy($ws.Range("Z1"))
$ws.Range("B3").Copy($ws.Range("Z2"))
$ws.Range("C2").Copy($ws.Range("Z3"))
$ws.Range("C3").Copy($ws.Range("Z4"))

# Remember C2/C3 formulas so their content can be restored exactly
# (Hyperlinks.Add with a TextToDisplay argument overwrites cell content).
$c2Formula = $ws.Range("C2").Formula
$c3Formula = $ws.Range("C3").Formula

$etsunan = "https://iiif.dl.itc.u-tokyo.ac.jp/iiif/tmp/toyo/suikei/Etsunan.tif/info.json"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), $etsunan) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), $etsunan) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), $etsunan, [Type]::Missing, [Type]::Missing, $etsunan) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), $etsunan, [Type]::Missing, [Type]::Missing, $etsunan) | Out-Null

$ws.Range("C2").Formula = $c2Formula
$ws.Range("C3").Formula = $c3Formula

# Restore the original hyperlink-cell formatting.
$ws.Range("Z1").Copy()
$ws.Range("B2").PasteSpecial(-4122) | Out-Null
$ws.Range("Z2").Copy()
$ws.Range("B3").PasteSpecial(-4122) | Out-Null
$ws.Range("Z3").Copy()
$ws.Range("C2").PasteSpecial(-4122) | Out-Null
$ws.Range("Z4").Copy()
$ws.Range("C3").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1:Z4").Clear()

# collection sheet is no longer the active tab; its stored selection moves to J14.
# (do this before activating "media" below, since selecting a range on a
# sheet implicitly makes that sheet active)
$ws5 = $wb.Worksheets.Item("collection")
$ws5.Range("J14").Select()

# media becomes the active sheet/tab, with B3 selected.
$ws.Activate()
$ws.Range("B3").Select()
